$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row at position 162, pushing existing rows 162..210 down to 163..211
$ws.Rows.Item(162).Insert()

# Populate the newly inserted row 162 with the new weekly record
$ws.Cells.Item(162, 1).Value = 4
$ws.Cells.Item(162, 2).Value = "Feria Lagunitas de Puerto Montt"
$ws.Cells.Item(162, 3).Value = "Los Lagos"
$ws.Cells.Item(162, 4).Value = 44588
$ws.Cells.Item(162, 5).Value = 10
$ws.Cells.Item(162, 6).Value = 100112003
$ws.Cells.Item(162, 7).Value = "Ajo"
$ws.Cells.Item(162, 8).Value = "Chino"
$ws.Cells.Item(162, 9).Value = "Primera"
$ws.Cells.Item(162, 10).Value = 80
$ws.Cells.Item(162, 11).Value = 21000
$ws.Cells.Item(162, 12).Value = 22000
$ws.Cells.Item(162, 13).Value = 21500
$ws.Cells.Item(162, 14).Value = "$/caja 10 kilos"
$ws.Cells.Item(162, 15).Value = "China"
$ws.Cells.Item(162, 16).Value = 2150
$ws.Cells.Item(162, 17).Value = 10
$ws.Cells.Item(162, 18).Value = "Hortaliza"
